# Generate Report for Handback
# Updates the "zh-cn" and "de-de" localization-status sheets with the
# handback information that has just become available for the
# "aaccdb49-5713-42ed-880c-d6d694af9c66" file (row 7 of each language sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# zh-cn sheet (row 7)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Latest Target File (I7) - becomes a hyperlink to the handed-back markdown file
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4f44abb3338592351cf8ff7a5c118c713d7fdc0c/e2e/aaccdb49-5713-42ed-880c-d6d694af9c66.md", "", "", "aaccdb49-5713-42ed-880c-d6d694af9c66.md")

# Latest Handback File (J7)
$wsZhCn.Range("J7").Value = "aaccdb49-5713-42ed-880c-d6d694af9c66.b0e70e3da30965a718a01fa6b4db791ba408eeb6.zh-cn.xlf"

# Latest Handback DateTime (K7)
$wsZhCn.Range("K7").Value = "2016-09-02 09:04:31"

# Error Detail (P7)
$wsZhCn.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f44abb3338592351cf8ff7a5c118c713d7fdc0c/e2e/aaccdb49-5713-42ed-880c-d6d694af9c66.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17dacca0a28517f8d7abe702b914413b0d61c861/e2e/aaccdb49-5713-42ed-880c-d6d694af9c66.md."

# ---------------------------------------------------------------------------
# de-de sheet (row 7)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest Target File (I7) - becomes a hyperlink to the handed-back markdown file
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4f44abb3338592351cf8ff7a5c118c713d7fdc0c/e2e/aaccdb49-5713-42ed-880c-d6d694af9c66.md", "", "", "aaccdb49-5713-42ed-880c-d6d694af9c66.md")

# Latest Handback File (J7)
$wsDeDe.Range("J7").Value = "aaccdb49-5713-42ed-880c-d6d694af9c66.b0e70e3da30965a718a01fa6b4db791ba408eeb6.de-de.xlf"

# Latest Handback DateTime (K7)
$wsDeDe.Range("K7").Value = "2016-09-02 09:04:39"

# Error Detail (P7) - same message as zh-cn (version mismatch against latest handoff)
$wsDeDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f44abb3338592351cf8ff7a5c118c713d7fdc0c/e2e/aaccdb49-5713-42ed-880c-d6d694af9c66.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17dacca0a28517f8d7abe702b914413b0d61c861/e2e/aaccdb49-5713-42ed-880c-d6d694af9c66.md."
